# 10/30 update for check
# Mark the newly-completed rows in the "二刷" (second pass) sheet by
# setting column D to 1 for the rows whose checkbox/flag cell was blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("二刷")

$ws.Range("D27").Value = 1
$ws.Range("D29").Value = 1
$ws.Range("D30").Value = 1
$ws.Range("D32").Value = 1
$ws.Range("D37").Value = 1
$ws.Range("D38").Value = 1
$ws.Range("D39").Value = 1
$ws.Range("D40").Value = 1

# Update the visible scroll position / active selection on that sheet to
# match where the user ended up after checking things off.
$ws.Activate()
$ws.Range("J34").Select()
$excel.ActiveWindow.ScrollRow = 25
